$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.541.48"
$ws.Range("E2").Value = "  -3.40%  "
$ws.Range("D3").Value = "1.781.18"
$ws.Range("E3").Value = "  -2.40%  "
$ws.Range("D4").Value = "'1.006"
$ws.Range("E4").Value = "  +0.31%  "
$ws.Range("D5").Value = "'1.005"
$ws.Range("E5").Value = "  +0.33%  "
$ws.Range("D6").Value = "'307.25"
$ws.Range("E6").Value = "  -1.76%  "
$ws.Range("D7").Value = "'0.4307"
$ws.Range("E7").Value = "  +1.50%  "
$ws.Range("D8").Value = "'0.3675"
$ws.Range("E8").Value = "  +2.39%  "
$ws.Range("D9").Value = "'0.07226"
$ws.Range("E9").Value = "  +0.40%  "
$ws.Range("D10").Value = "'0.8550"
$ws.Range("E10").Value = "  -0.29%  "
$ws.Range("D11").Value = "'20.46"
$ws.Range("E11").Value = "  -0.09%  "
$ws.Range("D12").Value = "1.773.77"
$ws.Range("E12").Value = "  -0.37%  "
$ws.Range("D13").Value = "'6.477"
$ws.Range("E13").Value = "  +0.31%  "
$ws.Range("D14").Value = "'5.275"
$ws.Range("E14").Value = "  -1.50%  "
$ws.Range("D15").Value = "'0.06917"
$ws.Range("E15").Value = "  -0.10%  "
$ws.Range("E16").Value = "  +0.40%  "
$ws.Range("D17").Value = "'79.66"
$ws.Range("E17").Value = "  -2.18%  "
$ws.Range("D18").Value = "'0.000008737"
$ws.Range("E18").Value = "  -1.83%  "
$ws.Range("E19").Value = "  +0.10%  "
$ws.Range("D20").Value = "'15.13"
$ws.Range("E20").Value = "  -1.53%  "
$ws.Range("D21").Value = "26.548.23"
$ws.Range("E21").Value = "  -2.68%  "
$ws.Range("D22").Value = "'5.128"
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("D23").Value = "'11.21"
$ws.Range("E23").Value = "  +2.90%  "
$ws.Range("D24").Value = "1.983.79"
$ws.Range("E24").Value = "  -1.85%  "
$ws.Range("D25").Value = "'152.55"
$ws.Range("E25").Value = "  -1.31%  "
$ws.Range("D26").Value = "'1.875"
$ws.Range("E26").Value = "  -5.84%  "
$ws.Range("D27").Value = "'18.20"
$ws.Range("E27").Value = "  -1.81%  "
$ws.Range("D28").Value = "'5.125"
$ws.Range("E28").Value = "  -0.28%  "
$ws.Range("D29").Value = "'114.61"
$ws.Range("E29").Value = "  +0.51%  "
$ws.Range("D30").Value = "'1.788"
$ws.Range("E30").Value = "  +0.50%  "
$ws.Range("D31").Value = "'0.08992"
$ws.Range("E31").Value = "  +0.90%  "
$ws.Range("D32").Value = "'0.7302"
$ws.Range("E32").Value = "  -1.94%  "
$ws.Range("E33").Value = "  +0.43%  "
$ws.Range("D34").Value = "'4.359"
$ws.Range("E34").Value = "  -3.85%  "
$ws.Range("D35").Value = "'2.754"
$ws.Range("E35").Value = "  -7.56%  "
$ws.Range("D37").Value = "'1.091"
$ws.Range("E37").Value = "  +1.02%  "
$ws.Range("D38").Value = "'0.05198"
$ws.Range("E38").Value = "  -0.94%  "
$ws.Range("D39").Value = "'0.01899"
$ws.Range("E39").Value = "  -1.37%  "
$ws.Range("D40").Value = "'0.4956"
$ws.Range("E40").Value = "  -2.04%  "
$ws.Range("D41").Value = "'0.1621"
$ws.Range("E41").Value = "  -2.13%  "
$ws.Range("D42").Value = "'2.607"
$ws.Range("E42").Value = "  -6.43%  "
$ws.Range("D43").Value = "'6.381"
$ws.Range("E43").Value = "  +0.73%  "
$ws.Range("D44").Value = "'8.096"
$ws.Range("E44").Value = "  -2.98%  "
$ws.Range("D45").Value = "'105.29"
$ws.Range("E45").Value = "  -0.99%  "
$ws.Range("D46").Value = "'10.25"
$ws.Range("E46").Value = "  -1.67%  "
$ws.Range("D47").Value = "'1.005"
$ws.Range("E47").Value = "  +0.41%  "
$ws.Range("D48").Value = "'0.4523"
$ws.Range("E48").Value = "  -3.13%  "
$ws.Range("D49").Value = "'0.06222"
$ws.Range("E49").Value = "  -3.62%  "
$ws.Range("D50").Value = "'1.617"
$ws.Range("E50").Value = "  +0.26%  "
$ws.Range("D51").Value = "'1.767"
$ws.Range("E51").Value = "  +4.63%  "
